# Auto-update the tracking sheet: decrement remaining days (column E) by 1
# for every data row. When a row's remaining count has hit 1 (i.e. the
# current cycle is about to roll over), reset it to a full new cycle of
# 10 remaining days and advance the start date (column F, stored as a
# plain YYYYMMDD integer) by 10 calendar days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header is row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)   # column D - 总天 (total days in cycle)
    $eCell = $ws.Cells.Item($r, 5)   # column E - 剩余 (remaining)
    $fCell = $ws.Cells.Item($r, 6)   # column F - 开始时间 (start date, YYYYMMDD)

    $eVal = $eCell.Value()
    if ($eVal -eq $null -or $eVal -eq "") { continue }

    $eVal = [int]$eVal

    $dVal = $dCell.Value()
    if ($dVal -ne $null -and $dVal -ne "" -and [int]$dVal -eq $eVal) {
        # Remaining equals the total cycle length: this row's countdown
        # has not started yet, so leave it untouched.
        continue
    }

    if ($eVal -eq 1) {
        # Cycle rolled over: reset remaining to a full new cycle and push
        # the start date forward by 10 days.
        $fVal = $fCell.Value()
        if ($fVal -ne $null -and $fVal -ne "") {
            $fVal = [int]$fVal
            $year = [int]([math]::Floor($fVal / 10000))
            $month = [int]([math]::Floor(($fVal % 10000) / 100))
            $day = [int]($fVal % 100)
            $date = Get-Date -Year $year -Month $month -Day $day
            $newDate = $date.AddDays(10)
            $newFVal = [int]($newDate.ToString("yyyyMMdd"))
            $fCell.Value = $newFVal
        }
        $eCell.Value = 10
    } else {
        $eCell.Value = $eVal - 1
    }
}
